$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.47123200046431
$ws.Range("C2").Value = 14.98209269522088
$ws.Range("E2").Value = 16.50058039820641
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.616019757334481
$ws.Range("I2").Value = 22.53030556941171
$ws.Range("N2").Value = 16.77605563675621
$ws.Range("B3").Value = 16.67661171436607
$ws.Range("C3").Value = 14.09459930748019
$ws.Range("E3").Value = 15.55601056631372
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.621022851344635
$ws.Range("I3").Value = 22.41029112157617
$ws.Range("N3").Value = 16.85416915921614
$ws.Range("B4").Value = 16.17490364407455
$ws.Range("C4").Value = 13.52498877252832
$ws.Range("E4").Value = 14.95233066869531
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.624241908211719
$ws.Range("I4").Value = 22.34477962614002
$ws.Range("N4").Value = 16.90415596332716
$ws.Range("B5").Value = 15.96730848111109
$ws.Range("C5").Value = 13.28687653984742
$ws.Range("E5").Value = 14.70063376233248
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.625590893867343
$ws.Range("I5").Value = 22.32014309014241
$ws.Range("N5").Value = 16.92503616239311
$ws.Range("B6").Value = 15.93265791655571
$ws.Range("C6").Value = 13.24698356458121
$ws.Range("E6").Value = 14.65850481261952
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.625817144021434
$ws.Range("I6").Value = 22.31617665501722
$ws.Range("N6").Value = 16.92853414694799
$ws.Range("B7").Value = 16.17211621335212
$ws.Range("C7").Value = 13.52180144654631
$ws.Range("E7").Value = 14.94895884913696
$ws.Range("F7").Value = 15.26647399323726
$ws.Range("G7").Value = 3.624259950247846
$ws.Range("I7").Value = 22.34443902671514
$ws.Range("N7").Value = 16.90443549393548
$ws.Range("B8").Value = 17.20029587808189
$ws.Range("C8").Value = 14.68132451820364
$ws.Range("E8").Value = 16.17995969068069
$ws.Range("F8").Value = 16.53996406344765
$ws.Range("G8").Value = 3.617714403877476
$ws.Range("I8").Value = 22.48723054006954
$ws.Range("N8").Value = 16.80256955054652
$ws.Range("B9").Value = 19.09472445415976
$ws.Range("C9").Value = 16.7518879634464
$ws.Range("E9").Value = 18.49033722141196
$ws.Range("F9").Value = 19.0027458068253
$ws.Range("G9").Value = 3.606037127601787
$ws.Range("I9").Value = 22.83180152499477
$ws.Range("N9").Value = 16.61882838767211
$ws.Range("B10").Value = 20.39853215296211
$ws.Range("C10").Value = 18.14172649495072
$ws.Range("E10").Value = 20.16125381810702
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.598151502938538
$ws.Range("I10").Value = 23.12362854162946
$ws.Range("N10").Value = 16.49352931191562
$ws.Range("B11").Value = 20.97032126172734
$ws.Range("C11").Value = 18.74442890357734
$ws.Range("E11").Value = 20.88015517379789
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.594712023132533
$ws.Range("I11").Value = 23.26456624362861
$ws.Range("N11").Value = 16.43861864165932
$ws.Range("B12").Value = 21.18362921188057
$ws.Range("C12").Value = 18.96834821035632
$ws.Range("E12").Value = 21.14650349627315
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.593430610018536
$ws.Range("I12").Value = 23.31908901199804
$ws.Range("N12").Value = 16.41812482395865
$ws.Range("B13").Value = 21.13783470375936
$ws.Range("C13").Value = 18.92031574462777
$ws.Range("E13").Value = 21.08940165247625
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.593705652556742
$ws.Range("I13").Value = 23.30729570599492
$ws.Range("N13").Value = 16.42252521668451
$ws.Range("B14").Value = 20.98793546646867
$ws.Range("C14").Value = 18.7629375813322
$ws.Range("E14").Value = 20.90218549736426
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.594606179919362
$ws.Range("I14").Value = 23.26902891106687
$ws.Range("N14").Value = 16.43692660090638
$ws.Range("B15").Value = 20.89569497100728
$ws.Range("C15").Value = 18.66597590351377
$ws.Range("E15").Value = 20.78674515519339
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.595160513414867
$ws.Range("I15").Value = 23.24573874572252
$ws.Range("N15").Value = 16.44578687052069
$ws.Range("B16").Value = 20.36072195204823
$ws.Range("C16").Value = 18.10173825351464
$ws.Range("E16").Value = 20.11344628155222
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.598379236662323
$ws.Range("I16").Value = 23.11458062022518
$ws.Range("N16").Value = 16.49715982152554
$ws.Range("B17").Value = 20.02695776282874
$ws.Range("C17").Value = 17.74798011609915
$ws.Range("E17").Value = 19.68987466591927
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.600391514381927
$ws.Range("I17").Value = 23.03619927796035
$ws.Range("N17").Value = 16.52920989091427
$ws.Range("B18").Value = 19.83298601304421
$ws.Range("C18").Value = 17.54172972540378
$ws.Range("E18").Value = 19.4423663988601
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.60156284044497
$ws.Range("I18").Value = 22.99188837559203
$ws.Range("N18").Value = 16.54784085991455
$ws.Range("B19").Value = 19.76697210798594
$ws.Range("C19").Value = 17.47142194486969
$ws.Range("E19").Value = 19.35789650991839
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.601961827679677
$ws.Range("I19").Value = 22.97701873923139
$ws.Range("N19").Value = 16.55418277598994
$ws.Range("B20").Value = 20.06269563150063
$ws.Range("C20").Value = 17.78592630121976
$ws.Range("E20").Value = 19.73536596341289
$ws.Range("F20").Value = 20.2495528364879
$ws.Range("G20").Value = 3.600175865059543
$ws.Range("I20").Value = 23.04446337501085
$ws.Range("N20").Value = 16.52577776221428
$ws.Range("B21").Value = 21.03205286925262
$ws.Range("C21").Value = 18.80928078319695
$ws.Range("E21").Value = 20.95733478031904
$ws.Range("F21").Value = 21.46857628470567
$ws.Range("G21").Value = 3.594341103698035
$ws.Range("I21").Value = 23.28023771870718
$ws.Range("N21").Value = 16.43268843754778
$ws.Range("B22").Value = 21.64678446842433
$ws.Range("C22").Value = 19.46554232198978
$ws.Range("E22").Value = 21.72168473016926
$ws.Range("F22").Value = 22.22866616901555
$ws.Range("G22").Value = 3.590650312042504
$ws.Range("I22").Value = 23.44103561952009
$ws.Range("N22").Value = 16.3735955693089
$ws.Range("B23").Value = 21.32045427464691
$ws.Range("C23").Value = 19.11173196272962
$ws.Range("E23").Value = 21.31685944640813
$ws.Range("F23").Value = 21.82633154475864
$ws.Range("G23").Value = 3.592609009071165
$ws.Range("I23").Value = 23.35460989344616
$ws.Range("N23").Value = 16.40497495570403
$ws.Range("B24").Value = 20.04654502767191
$ws.Range("C24").Value = 17.76877976578864
$ws.Range("E24").Value = 19.7148117829744
$ws.Range("F24").Value = 20.22900810905294
$ws.Range("G24").Value = 3.600273315123081
$ws.Range("I24").Value = 23.04072483394773
$ws.Range("N24").Value = 16.52732878915311
$ws.Range("B25").Value = 18.59683088616108
$ws.Range("C25").Value = 16.21439311147944
$ws.Range("E25").Value = 17.8381906676167
$ws.Range("F25").Value = 18.34778573295697
$ws.Range("G25").Value = 3.609073410736226
$ws.Range("I25").Value = 22.73173117267301
$ws.Range("N25").Value = 16.66682776032881
